# Update crypto price/volume data to the latest scrape
# (Sat Dec 24 17:14:16 UTC 2022).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to keep storing numeric-looking strings as text
    # (matching how the original sheet stores prices), then drop the
    # quote-prefix/number-format styling Excel would otherwise apply
    # so the cell keeps using the default style.
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# Simple price/volume refreshes (no row reshuffle)
Set-TextValue "D2" "244.58"
Set-TextValue "D3" "21.85"
Set-TextValue "D4" "5.389"
Set-TextValue "D7" "0.8152"
Set-TextValue "D8" "0.9486"

# Rows 9-17: the coin list rotated by one position (One moved from row 9 to
# row 17, and every other coin shifted up by one row), each with refreshed
# price/volume data.
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D9" "0.1437"
$ws.Range("E9").Value = "8WazirXWRX"

$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D10" "0.07405"
$ws.Range("E10").Value = "9MandalaExchangeTokenMDX"

$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D11" "0.03467"
$ws.Range("E11").Value = "10LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D12" "0.03055"
$ws.Range("E12").Value = "11BitrueCoinBTR"

$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D13" "0.09412"
$ws.Range("E13").Value = "12BitMartTokenBMX"

$ws.Range("B14").Value = "MCDex"
$ws.Range("C14").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D14" "4.005"
$ws.Range("E14").Value = "13MCDexMCB"

$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D15" "0.001592"
$ws.Range("E15").Value = "14BitForexTokenBF"

$ws.Range("B16").Value = "CoinExToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D16" "0.04800"
$ws.Range("E16").Value = "15CoinExTokenCET"

$ws.Range("B17").Value = "One"
$ws.Range("C17").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D17" "0.0005951"
$ws.Range("E17").Value = "16OneONE"

# Remaining simple price refreshes further down the sheet
Set-TextValue "D18" "0.005594"
Set-TextValue "D19" "0.004160"
Set-TextValue "D20" "0.0009864"
Set-TextValue "D21" "3.672"
Set-TextValue "D22" "6.423"
Set-TextValue "D23" "2.174"
Set-TextValue "D26" "0.00007001"
Set-TextValue "D40" "0.04014"
Set-TextValue "D41" "0.006401"
Set-TextValue "D43" "0.002900"
Set-TextValue "D44" "0.006649"
$ws.Range("E44").Value = "43LocalTradersLCTBestin24h"
Set-TextValue "D45" "0.00005240"
Set-TextValue "D48" "0.002788"
Set-TextValue "D49" "0.00002100"
